# Auto-generated edit script: updates Price (D) and Volume(1h) (E) columns
# per the commit diff ("Updated cryptos list ... with GitHub Actions").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.629.95'
$ws.Range("E2").Value = '  +0.16%  '

$ws.Range("D3").Value = '2.649.18'
$ws.Range("E3").Value = '  -0.51%  '

$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '596.55'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.32%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '156.39'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +0.07%  '

$ws.Range("E7").Value = '  +0.01%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.626'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +4.08%  '

$ws.Range("E9").Value = '  +4.46%  '

$ws.Range("E10").Value = '  +0.44%  '

$ws.Range("E11").Value = '  -1.35%  '

$ws.Range("E12").Value = '  +1.02%  '

$ws.Range("E13").Value = '  -1.46%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000196'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +1.25%  '

$ws.Range("D15").Value = '3.128.58'
$ws.Range("E15").Value = '  -0.37%  '

$ws.Range("D16").Value = '65.497.75'
$ws.Range("E16").Value = '  +0.21%  '

$ws.Range("D17").Value = '2.656.78'
$ws.Range("E17").Value = '  -0.07%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.59'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +1.12%  '

$ws.Range("E19").Value = '  -0.98%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.46'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -0.04%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '349.27'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -0.04%  '

$ws.Range("E22").Value = '  +0.14%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '69.16'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -0.67%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.0000112'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +5.52%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.68'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +0.28%  '

$ws.Range("E26").Value = '  +1.43%  '

$ws.Range("E27").Value = '  +0.11%  '

$ws.Range("E28").Value = '  -1.44%  '

$ws.Range("E29").Value = '  +0.11%  '

$ws.Range("E30").Value = '  -1.69%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '529.78'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -1.84%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.12'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -1.67%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.75'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +0.74%  '

$ws.Range("E34").Value = '  -1.08%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.43'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +0.64%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.419'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -0.13%  '

$ws.Range("E37").Value = '  +0.23%  '

$ws.Range("E38").Value = '  +0.04%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.92'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -0.51%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '153.43'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -3.42%  '

$ws.Range("E41").Value = '  +0.08%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '160.40'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -3.09%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.06'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +0.16%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.29'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +1.99%  '

$ws.Range("E45").Value = '  -0.53%  '

$ws.Range("E46").Value = '  -1.40%  '

$ws.Range("E47").Value = '  -1.71%  '

$ws.Range("E48").Value = '  -1.64%  '

$ws.Range("E49").Value = '  +12.40%  '

$ws.Range("E50").Value = '  +0.03%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '19.71'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -0.74%  '
